$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH_WF_REGISTRATION_REVIEW")
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# ---------------------------------------------------------------------------
# Sheet 1: LH_WF_REGISTRATION_REVIEW
# ---------------------------------------------------------------------------

# Row 2: the review ID moves to the new dash-separated naming convention;
# the rest of the row keeps its existing content.
$ws1.Range("B2").Value = "LH-WF-REGISTRATION-REVIEW-001"

# Row 3: a second review round was added (previously a blank placeholder
# row). Pull the cell formatting from row 2 so fills/borders match, then
# overwrite the B3 cell (whose fill differs slightly from the rest of the
# row) with the same style as B2.
$ws1.Range("B2").Copy()
$ws1.Range("B3").PasteSpecial(-4122)

$ws1.Range("A3").Value = "27/4/2025"
$ws1.Range("B3").Value = "LH-WF-REGISTRATION-REVIEW-002"
$ws1.Range("C3").Value = "Ahmed Abuzaid"
$ws1.Range("D3").Value = "v1.1"
$ws1.Range("E3").Value = "back to SRS I found there are many error messages for many validations like existing user name or existing email or validations fro password"
$ws1.Range("F3").Value = "so I prefer to add some error messages ""with red color"" from the SRS to wireframe to make it more expressive, you can back to login wireframe to understand what I mean"
$ws1.Range("G3").Value = "Gehad"
$ws1.Range("H3").Value = "closed"
$ws1.Range("I3").Value = "closed"

$ws1.Rows.Item(3).RowHeight = 150

# Leave the sheet scrolled back to column A and with C3 as the last
# selected cell (matches the saved view state).
$ws1.Range("C3").Select()

# ---------------------------------------------------------------------------
# Sheet 2: VERSION-HISTORY
# ---------------------------------------------------------------------------

# Rows 4 and 5: two more version-history entries were added. Copy the
# formatting from row 2 (a fully styled existing row, using the plain
# "Updated section" fill rather than row 3's alternate-shading one) down
# into the new rows before writing their values.
$ws2.Range("A2:D2").Copy()
$ws2.Range("A4:D4").PasteSpecial(-4122)
$ws2.Range("A5:D5").PasteSpecial(-4122)

$ws2.Range("A4").Value = "v1.2"
$ws2.Range("B4").Value = "Ahmed Abuzaid"
$ws2.Range("C4").Value = "ask to add more details to registration form wireframe"
$ws2.Range("D4").Value = 45774

$ws2.Range("A5").Value = "v1.3"
$ws2.Range("B5").Value = "Ahmed Abuzaid"
$ws2.Range("C5").Value = "close registration wireframe review, verify the updates and modify id naming convention"
$ws2.Range("D5").Value = 45775

$ws2.Rows.Item(4).RowHeight = 37.5
$ws2.Rows.Item(5).RowHeight = 56.25

# VERSION-HISTORY stays the active sheet/tab; select the new rows.
$ws2.Range("C12:C13").Select()
